$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update row 14: "Applications web State-of-the-Art, Architecture" moves from
#    "En cours" to "Termine", with realised hours 0.5 -> 3.
$ws.Range("C14").Value = "Terminé"
$ws.Range("E14").Value = 3

# 2. Insert a new row at 15 for the new task "Charte graphique" and give it the
#    same look (borders/number format) as its neighbours (A14 and F16).
$ws.Rows.Item(15).Insert()
$ws.Range("F16").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)

$ws.Range("A15").Value = "Analyse"
$ws.Range("B15").Value = "Charte graphique"
$ws.Range("C15").Value = "En cours"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 1
$ws.Range("F15").Formula = "=E15/D15"

# 3. Old row 16 ("Frameworks JavaScript, Front-end") is now row 17: realised
#    hours 0.5 -> 1.5.
$ws.Range("E17").Value = 1.5

# 4. Old row 19 (blank task, domaine "Analyse") is now row 20: planned hours
#    87 -> 80.
$ws.Range("D20").Value = 80
